# 202101_biden_cabinet_confirmation_predictions_observed.xlsx
# Adds observed confirmation dates (column I) + their computed error in days
# (column J) for several cabinet positions, and adds a new "observed_error"
# summary sheet (mean error / rmse / mae over column J).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)      # 202101_biden_cabinet_confirmati
$ws2 = $wb.Worksheets.Item(2)      # model_errors_by_position

# ---------------------------------------------------------------------
# 1. Fill in newly-observed confirmation dates on the main sheet and the
#    error-in-days formula that depends on them (=I<r> - F<r>).
# ---------------------------------------------------------------------
$observed = @{
    2  = 44229
    7  = 44235
    9  = 44250
    10 = 44229
    12 = 44252
    13 = 44256
    14 = 44257
}

foreach ($r in $observed.Keys) {
    $ws1.Range("I$r").Value = $observed[$r]
    $ws1.Range("J$r").Formula = "=I$r - F$r"
}

# J2 inherits the date number-format from before (it used to be an empty
# cell formatted like the date columns); once it holds the day-count
# formula it should look like the rest of the (already-populated) J
# column, e.g. J3 -- so copy that cell's formatting over.
$ws1.Range("J3").Copy() | Out-Null
$ws1.Range("J2").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# Move the selection like the author left it.
$ws1.Range("A11").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. Add the new "observed_error" summary sheet (after the existing two).
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "observed_error"

$ws3.Range("A1").Value = "mean_error"
$ws3.Range("B1").Value = "rmse"
$ws3.Range("C1").Value = "mae"

$ws3.Range("A2").Formula = "=AVERAGE('202101_biden_cabinet_confirmati'!J2:J16)"
$ws3.Range("A2").ClearFormats() | Out-Null

$ws3.Range("B2").Formula = "=SQRT(AVERAGE(SUMSQ('202101_biden_cabinet_confirmati'!J2:J16)))"

$ws3.Range("C2").FormulaArray = "=AVERAGE(ABS('202101_biden_cabinet_confirmati'!J2:J16))"

$ws3.Range("E18").Select() | Out-Null

# Leave the original sheet as the active / visible tab again.
$ws1.Select() | Out-Null
